# Update the "Latest/Correspond Xliff Generate/Handoff/Handback" timestamp
# cells on the handback-status report, as produced by a fresh report
# generation run ("Generate Report for Handback").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 0e0c301f-...md
$wsOverview.Range("G2").Value = "2016-08-15 11:02:44"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the same file
$wsZhCn.Range("H2").Value = "2016-08-15 11:02:39"
$wsZhCn.Range("K2").Value = "2016-08-15 11:02:55"

# de-de sheet: Correspond Handoff / Handback datetimes for the same file
$wsDeDe.Range("H2").Value = "2016-08-15 11:02:44"
$wsDeDe.Range("K2").Value = "2016-08-15 11:03:07"
